# LOB1056.docx restructuring script
# The document's "Heading2" section headers (Objetivos, Docente(s) Responsavel(eis),
# Programa resumido, Programa, Avaliacao, Bibliografia, Requisitos) all stay in their
# original paragraph slots. Only the body/content paragraphs between those headings
# are rearranged - their text (and, for one paragraph, internal bold runs) is swapped
# around, while each paragraph keeps its own original paragraph style.

$d = $word.ActiveDocument

function Set-ParaText {
    # Replaces the full text of paragraph number $ParaIndex with $NewText.
    # Embedded "<BR>" markers in $NewText become manual line breaks (Chr(11), i.e. <w:br/>).
    param($Doc, $ParaIndex, $NewText)

    $p = $Doc.Paragraphs.Item($ParaIndex)
    $pr = $p.Range
    $startPos = $pr.Start
    $endPos = $pr.End

    $full = $NewText -replace "<BR>", [string][char]11

    $rng = $Doc.Range($startPos, $endPos)
    $rng.Text = $full
}

function Set-ParaSegments {
    # Replaces the full text of paragraph number $ParaIndex with the concatenation of
    # $Segments (array of @{Text=...; Bold=$true/$false}), applying Bold explicitly to
    # each segment so the correct runs end up bold / not bold. "<BR>" -> manual line break.
    param($Doc, $ParaIndex, $Segments)

    $p = $Doc.Paragraphs.Item($ParaIndex)
    $pr = $p.Range
    $startPos = $pr.Start
    $endPos = $pr.End

    $full = ""
    foreach ($seg in $Segments) {
        $full += ($seg.Text -replace "<BR>", [string][char]11)
    }

    $rng = $Doc.Range($startPos, $endPos)
    $rng.Text = $full

    $pos = $startPos
    foreach ($seg in $Segments) {
        $t = $seg.Text -replace "<BR>", [string][char]11
        $len = $t.Length
        if ($len -gt 0) {
            $segRng = $Doc.Range($pos, $pos + $len)
            if ($seg.Bold) {
                $segRng.Font.Bold = $true
            } else {
                $segRng.Font.Bold = $false
            }
        }
        $pos += $len
    }
}

# --- Paragraph 6 (under "Objetivos"): now holds the PT "Programa resumido" text ---
$txt6 = "Descrição do programa resumido em português.<BR>" + `
  "Introdução ao Matlab (ou SciLab, Octave, Freemat, etc.), raízes de equações não lineares, sistemas de equações, ajuste de curvas pelo método dos mínimos quadrados, integração numérica e equações diferenciais ordinárias."
Set-ParaText $d 6 $txt6

# --- Paragraph 7 (under "Objetivos", italic): now holds the EN "Programa resumido" text ---
$txt7 = "Introduction to Matlab (or SciLab, Octave, Freemat, etc.), roots of nonlinear equations, systems of equations, least-squares fitting of curves to data, numerical integration, and solving ordinary differential equations."
Set-ParaText $d 7 $txt7

# --- Paragraph 9 (under "Docente(s) Responsável(eis)", ListBullet): now holds the PT "Objetivos" text ---
$txt9 = "Descrição dos objetivos em português.<BR>" + `
  "Aplicar os princípios e a lógica de programação de computadores no desenvolvimento de códigos para cálculo e visualização de dados. Usar ferramentas computacionais para resolver problemas em ciências e engenharia envolvendo: raízes de equações não lineares, sistemas de equações, ajuste de curvas pelo método dos mínimos quadrados, integração numérica e equações diferenciais ordinárias. Editores de texto e planilhas eletrônicas serão usados, fora do horário de aula, na preparação de relatórios técnicos, pôsteres, apresentações gráficas, etc."
Set-ParaText $d 9 $txt9

# --- Paragraph 11 (under "Programa resumido"): now holds the PT "Programa" (full) text ---
$txt11 = "Descrição do programa em português.<BR>" + `
  "1.Introdução ao Matlab (ou SciLab, Octave, Freemat, etc.); Cálculos simples; Uso de variáveis e funções nativas. Matrizes e Vetores; Gráficos.<BR>" + `
  "2.Programação em Matlab; escrevendo programas e funções, comandos de entrada e saída de dados, controle de fluxo de execução, vetorização e variáveis globais.<BR>" + `
  "3.Raízes de equações não-Lineares: ponto fixo, bissecção, e método de Newton.<BR>" + `
  "4.Álgebra Linear (propriedades e operações com matrizes e vetores).<BR>" + `
  "5.Resolução de sistemas de equações lineares.<BR>" + `
  "6.Ajuste de curvas pelo método dos mínimos quadrados.<BR>" + `
  "7.Interpolação.<BR>" + `
  "8.Integração numérica.<BR>" + `
  "9.Resolução de equações diferenciais ordinárias."
Set-ParaText $d 11 $txt11

# --- Paragraph 12 (under "Programa resumido", italic): now holds the EN "Objetivos" text ---
$txt12 = "Apply computer programming logic and principles to write code to perform computations and data visualization to solve a variety of science and engineering problems involving: roots of nonlinear equations, systems of equations, least-squares fitting of curves to data, numerical integration, and ordinary differential equations. Computational tools such as word processors and spreadsheets will be used outside of class to prepare technical reports, posters, presentations, etc."
Set-ParaText $d 12 $txt12

# --- Paragraph 14 (under "Programa"): now holds the "Método" text moved out of "Avaliação" ---
$txt14 = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
Set-ParaText $d 14 $txt14

# --- Paragraph 17 (under "Avaliação", ListBullet, mixed bold): shuffled Método/Critério/Norma
#     content, now also absorbing the "Bibliografia" reference list text ---
$txt17bib = "CHAPRA, Steven C., Métodos numéricos aplicados com MATLAB para engenheiros e cientistas. Porto Alegre: AMGH, 2013.<BR>" + `
               "CHAPRA, Steven C.; CANALE, Raymond P., Métodos numéricos para engenharia. São Paulo: McGraw-Hill, 2008.<BR>" + `
               "SPERANDIO, Decio; MENDES, João Teixeira, SILVA, Luiz Henry Monken. Cálculo numérico. São Paulo: Prentice Hall Brasil, 2003.<BR>" + `
               "RUGGIERO, M.A.G.; DA ROCHA LOPES, V.L. Cálculo Numérico: Aspectos Teóricos e Computacionais. Makron Books, 2a Edição, 1997.<BR>" + `
               "FRANCO, Neide Maria Bertoldi. Cálculo numérico. Prentice Hall Brasil, 2006.<BR>" + `
               "BURIAN, Reinaldo; LIMA, Antonio Carlos. Cálculo numérico. São Paulo: LTC, 2007.<BR>" + `
               "ARENALES, Selma Helena De Vasconcelos; DAREZZO, Artur. Cálculo Numérico. São Paulo: Thomson Pioneira, 2007."
$segs17 = @(
    @{ Text = "Método: "; Bold = $true },
    @{ Text = "NF≥ 5,0.<BR>"; Bold = $false },
    @{ Text = "Critério: "; Bold = $true },
    @{ Text = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.<BR>"; Bold = $false },
    @{ Text = "Norma de recuperação: "; Bold = $true },
    @{ Text = $txt17bib; Bold = $false }
)
Set-ParaSegments $d 17 $segs17

# --- Paragraph 19 (under "Bibliografia"): now holds the teacher name moved out of "Docente(s)" ---
$txt19 = "8870322 - Fabiano Fernandes Bargos"
Set-ParaText $d 19 $txt19
